$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# fix(FN-3460): fix invalid facility utilisation values -- all need to match
# as same facility id for all rows
$ws.Range("E5").Value = 600000
$ws.Range("G5").Value = 3938753.8
$ws.Range("E6").Value = 600000
$ws.Range("G6").Value = 761579.37

# Facility utilisation column (G) now holds values of the same order of
# magnitude as columns E/F, so widen it to line up with them
$ws.Columns("G:G").ColumnWidth = 15.5

# Reflect the cell range that was being reviewed/edited when the fix was made
$ws.Range("E5:H6").Select()
